# Rename the lone (empty) worksheet to "no_data" - it stays empty, it just
# demonstrates the "no data at all" case for read_excel / read_ods.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "no_data"

# Add a second worksheet, placed right after "no_data", that has a header
# row (column names) but zero data rows - the "no rows" case.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "no_rows"

$ws2.Range("A1").Value = "colx"
$ws2.Range("B1").Value = "coly"
$ws2.Range("C1").Value = "colz"

# Header row is bold, like a normal table header.
$ws2.Range("A1:C1").Font.Bold = $true
